$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 43 (hunk 0)
$ws.Range("H43").Value = 6280
$ws.Range("I43").Value = 2200
$ws.Range("J43").Value = 9000
$ws.Range("K43").Value = 2200
$ws.Range("L43").Value = 9000
$ws.Range("M43").Value = -2131
$ws.Range("N43").Value = -9138

# row 62 (hunk 1)
$ws.Range("H62").Value = 23887.092
$ws.Range("I62").Value = 20665.834
$ws.Range("K62").Value = 20665.834
$ws.Range("M62").Value = -20041.834

# row 65 (hunk 2)
$ws.Range("H65").Value = 23887.092
$ws.Range("I65").Value = 20665.834
$ws.Range("K65").Value = 103329.17
$ws.Range("M65").Value = -100209.17

# row 86 (hunk 3)
$ws.Range("H86").Value = 2231.2
$ws.Range("I86").Value = 1789
$ws.Range("K86").Value = 1789
$ws.Range("M86").Value = -666

# row 88 (hunk 4)
$ws.Range("H88").Value = 204697.8
$ws.Range("I88").Value = 1750
$ws.Range("J88").Value = 339996.34
$ws.Range("K88").Value = 1750
$ws.Range("L88").Value = 339996.34
$ws.Range("M88").Value = -1344
$ws.Range("N88").Value = -340808.34

# row 89 (hunk 5)
$ws.Range("H89").Value = 2231.2
$ws.Range("I89").Value = 1789
$ws.Range("K89").Value = 8945
$ws.Range("M89").Value = -3329

# row 91 (hunk 6)
$ws.Range("H91").Value = 204697.8
$ws.Range("I91").Value = 1750
$ws.Range("J91").Value = 339996.34
$ws.Range("K91").Value = 1750
$ws.Range("L91").Value = 339996.34
$ws.Range("M91").Value = -346
$ws.Range("N91").Value = -342804.34

# row 138 (hunk 7)
$ws.Range("H138").Value = 2041.9333
$ws.Range("J138").Value = 4096
$ws.Range("L138").Value = 12288
$ws.Range("N138").Value = -22568

$ws = $wb.Worksheets.Item("ARM")
# row 31 (hunk 8)
$ws.Range("H31").Value = 6874.4287
$ws.Range("I31").Value = 2933.6667
$ws.Range("J31").Value = 30519
$ws.Range("K31").Value = 2933.6667
$ws.Range("L31").Value = 30519
$ws.Range("M31").Value = -2639.6667
$ws.Range("N31").Value = -31107

# row 32 (hunk 9)
$ws.Range("H32").Value = 9721128
$ws.Range("I32").Value = 1726322.4
$ws.Range("J32").Value = 29473000
$ws.Range("K32").Value = 1726322.4
$ws.Range("L32").Value = 29473000
$ws.Range("M32").Value = -1726035.4
$ws.Range("N32").Value = -29473574

# row 61 (hunk 10)
$ws.Range("H61").Value = 2524.5833
$ws.Range("I61").Value = 2447.2
$ws.Range("J61").Value = 2911.5
$ws.Range("K61").Value = 2447.2
$ws.Range("L61").Value = 2911.5
$ws.Range("M61").Value = -2235.2
$ws.Range("N61").Value = -3335.5

# row 103 (hunk 11)
$ws.Range("H103").Value = 30180.5
$ws.Range("J103").Value = 30180.5
$ws.Range("L103").Value = 30180.5
$ws.Range("N103").Value = -32524.5

# row 110 (hunk 12)
$ws.Range("H110").Value = 1369.742
$ws.Range("I110").Value = 1035.36
$ws.Range("K110").Value = 1035.36
$ws.Range("M110").Value = 1009.64

# row 119 (hunk 13)
$ws.Range("H119").Value = 45098.668
$ws.Range("J119").Value = 45098.668
$ws.Range("L119").Value = 45098.668
$ws.Range("N119").Value = -54774.668

# row 132 (hunk 14)
$ws.Range("H132").Value = 4272.2905
$ws.Range("I132").Value = 4182
$ws.Range("J132").Value = 4581.857
$ws.Range("K132").Value = 12546
$ws.Range("L132").Value = 13745.571
$ws.Range("M132").Value = -10016
$ws.Range("N132").Value = -18805.571

# row 136 (hunk 15)
$ws.Range("H136").Value = 2524.5833
$ws.Range("I136").Value = 2447.2
$ws.Range("J136").Value = 2911.5
$ws.Range("K136").Value = 7341.599999999999
$ws.Range("L136").Value = 8734.5
$ws.Range("M136").Value = -4791.599999999999
$ws.Range("N136").Value = -13834.5

$ws = $wb.Worksheets.Item("BSM")
# row 20 (hunk 16)
$ws.Range("H20").Value = 3864.9773
$ws.Range("J20").Value = 4391.5
$ws.Range("L20").Value = 4391.5
$ws.Range("N20").Value = -4885.5

# row 22 (hunk 17)
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -327

# row 86 (hunk 18)
$ws.Range("H86").Value = 571.41174
$ws.Range("I86").Value = 433.8
$ws.Range("J86").Value = 1603.5
$ws.Range("K86").Value = 433.8
$ws.Range("L86").Value = 1603.5
$ws.Range("M86").Value = 689.2
$ws.Range("N86").Value = -3849.5

# row 89 (hunk 19)
$ws.Range("H89").Value = 571.41174
$ws.Range("I89").Value = 433.8
$ws.Range("J89").Value = 1603.5
$ws.Range("K89").Value = 2169
$ws.Range("L89").Value = 8017.5
$ws.Range("M89").Value = 3447
$ws.Range("N89").Value = -19249.5

# row 94 (hunk 20)
$ws.Range("H94").Value = 1333.1951
$ws.Range("I94").Value = 1284.742
$ws.Range("J94").Value = 1483.4
$ws.Range("K94").Value = 1284.742
$ws.Range("L94").Value = 1483.4
$ws.Range("M94").Value = -833.742
$ws.Range("N94").Value = -2385.4

# row 134 (hunk 21)
$ws.Range("H134").Value = 6765.25
$ws.Range("I134").Value = 5681.85
$ws.Range("J134").Value = 9473.75
$ws.Range("K134").Value = 17045.55
$ws.Range("L134").Value = 28421.25
$ws.Range("M134").Value = -14510.55
$ws.Range("N134").Value = -33491.25

$ws = $wb.Worksheets.Item("CRP")
# row 58 (hunk 22)
$ws.Range("H58").Value = 1212.8572
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 1372.5
$ws.Range("K58").Value = 1000
$ws.Range("L58").Value = 1372.5
$ws.Range("M58").Value = -797
$ws.Range("N58").Value = -1778.5

# row 132 (hunk 23)
$ws.Range("H132").Value = 3732.3333
$ws.Range("I132").Value = 3826.2
$ws.Range("K132").Value = 11478.6
$ws.Range("M132").Value = -8948.599999999999

# row 136 (hunk 24)
$ws.Range("H136").Value = 1212.8572
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 1372.5
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 4117.5
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -9217.5

$ws = $wb.Worksheets.Item("CUL")
# row 37 (hunk 25)
$ws.Range("H37").Value = 125085624
$ws.Range("J37").Value = 125085624
$ws.Range("L37").Value = 375256872
$ws.Range("N37").Value = -375257096

# row 109 (hunk 26)
$ws.Range("H109").Value = 59306.06
$ws.Range("I109").Value = 59306.06
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 177918.18
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -176878.18
$ws.Range("N109").ClearContents()

# row 134 (hunk 27)
$ws.Range("H134").Value = 2343.1667
$ws.Range("I134").Value = 2343.1667
$ws.Range("K134").Value = 7029.500100000001
$ws.Range("M134").Value = -1959.500100000001

# row 136 (hunk 28)
$ws.Range("H136").Value = 3396.2222
$ws.Range("I136").Value = 2080.8572
$ws.Range("K136").Value = 6242.571599999999
$ws.Range("M136").Value = -1142.571599999999

# row 138 (hunk 29)
$ws.Range("H138").Value = 1450.0834
$ws.Range("I138").Value = 1180.2
$ws.Range("K138").Value = 3540.6
$ws.Range("M138").Value = 1599.4

$ws = $wb.Worksheets.Item("GSM")
# row 2 (hunk 30)
$ws.Range("H2").Value = 431.58823
$ws.Range("I2").Value = 404
$ws.Range("J2").Value = 471
$ws.Range("K2").Value = 404
$ws.Range("L2").Value = 471
$ws.Range("M2").Value = -291
$ws.Range("N2").Value = -697

# row 113 (hunk 31)
$ws.Range("H113").Value = 7782.6313
$ws.Range("I113").Value = 9692
$ws.Range("K113").Value = 9692
$ws.Range("M113").Value = -7522

# row 132 (hunk 32)
$ws.Range("H132").Value = 2905.2273
$ws.Range("I132").Value = 2898.2363
$ws.Range("K132").Value = 8694.7089
$ws.Range("M132").Value = -6164.7089

$ws = $wb.Worksheets.Item("LTW")
# row 7 (hunk 33)
$ws.Range("H7").Value = 34836.125
$ws.Range("I7").Value = 38812.715
$ws.Range("K7").Value = 38812.715
$ws.Range("M7").Value = -38700.715

# row 16 (hunk 34)
$ws.Range("H16").Value = 2870
$ws.Range("I16").Value = 3114.4285
$ws.Range("K16").Value = 3114.4285
$ws.Range("M16").Value = -2944.4285

# row 68 (hunk 35)
$ws.Range("H68").Value = 39189.5
$ws.Range("J68").Value = 75250
$ws.Range("L68").Value = 75250
$ws.Range("N68").Value = -76748

# row 71 (hunk 36)
$ws.Range("H71").Value = 39189.5
$ws.Range("J71").Value = 75250
$ws.Range("L71").Value = 376250
$ws.Range("N71").Value = -383738

# row 82 (hunk 37)
$ws.Range("H82").Value = 815.8
$ws.Range("I82").Value = 610.5833
$ws.Range("K82").Value = 610.5833
$ws.Range("M82").Value = -249.5833

# row 85 (hunk 38)
$ws.Range("H85").Value = 815.8
$ws.Range("I85").Value = 610.5833
$ws.Range("K85").Value = 610.5833
$ws.Range("M85").Value = 637.4167

# row 93 (hunk 39)
$ws.Range("H93").Value = 1550.65
$ws.Range("I93").Value = 942.2308
$ws.Range("K93").Value = 942.2308
$ws.Range("M93").Value = 305.7692

# row 100 (hunk 40)
$ws.Range("H100").Value = 21032
$ws.Range("I100").Value = 14508.739
$ws.Range("J100").Value = 58540.75
$ws.Range("K100").Value = 14508.739
$ws.Range("L100").Value = 58540.75
$ws.Range("M100").Value = -13967.739
$ws.Range("N100").Value = -59622.75

# row 122 (hunk 41)
$ws.Range("H122").Value = 22666.666
$ws.Range("I122").Value = 52250
$ws.Range("J122").Value = 7875
$ws.Range("K122").Value = 156750
$ws.Range("L122").Value = 23625
$ws.Range("M122").Value = -154300
$ws.Range("N122").Value = -28525

# row 126 (hunk 42)
$ws.Range("H126").Value = 34836.125
$ws.Range("I126").Value = 38812.715
$ws.Range("K126").Value = 116438.145
$ws.Range("M126").Value = -113968.145

# row 132 (hunk 43)
$ws.Range("H132").Value = 4761.3105
$ws.Range("I132").Value = 4366.381
$ws.Range("J132").Value = 5798
$ws.Range("K132").Value = 13099.143
$ws.Range("L132").Value = 17394
$ws.Range("M132").Value = -10569.143
$ws.Range("N132").Value = -22454

$ws = $wb.Worksheets.Item("WVR")
# row 62 (hunk 44)
$ws.Range("H62").Value = 87437.336
$ws.Range("I62").Value = 15428.143
$ws.Range("J62").Value = 123441.93
$ws.Range("K62").Value = 15428.143
$ws.Range("L62").Value = 123441.93
$ws.Range("M62").Value = -14804.143
$ws.Range("N62").Value = -124689.93

# row 65 (hunk 45)
$ws.Range("H65").Value = 87437.336
$ws.Range("I65").Value = 15428.143
$ws.Range("J65").Value = 123441.93
$ws.Range("K65").Value = 77140.715
$ws.Range("L65").Value = 617209.6499999999
$ws.Range("M65").Value = -74020.715
$ws.Range("N65").Value = -623449.6499999999

# row 81 (hunk 46)
$ws.Range("H81").Value = 2500500
$ws.Range("I81").Value = 2500500
$ws.Range("K81").Value = 5001000
$ws.Range("M81").Value = -4999939

# row 84 (hunk 47)
$ws.Range("H84").Value = 2500500
$ws.Range("I84").Value = 2500500
$ws.Range("K84").Value = 25005000
$ws.Range("M84").Value = -24999696

# row 132 (hunk 48)
$ws.Range("H132").Value = 8636.522999999999
$ws.Range("I132").Value = 13162.154
$ws.Range("K132").Value = 39486.462
$ws.Range("M132").Value = -36956.462
